# Updates cryptos list figures (prices / 1h volume %) and fixes the
# Maker / RenderToken row ordering, matching the GitHub Actions data refresh.
# For "Price" (column D) cells whose new value looks like a plain number,
# force a Text number format before assigning so Excel doesn't silently
# coerce the string into a float (which would lose formatting such as
# "1.402.68" grouping or trailing zeros like "0.0530"). ClearFormats()
# afterwards drops the temporary format so the cell keeps the workbook's
# default (unstyled) appearance, same as the original cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.730.49'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.819.31'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.29'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.579'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.02%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '34.94'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +7.20%  '
$ws.Range('E9').Value = '  +1.63%  '
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0953'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.084.29'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.38'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.807.93'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.646'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.712.62'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('E17').Value = '  +1.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.28'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0801'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '246.63'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.59'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.17%  '
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '173.33'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +5.85%  '
$ws.Range('E25').Value = '  +2.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.53'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.85'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.79%  '
$ws.Range('E28').Value = '  +2.60%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  +2.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0530'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.66%  '
$ws.Range('E32').Value = '  +1.16%  '
$ws.Range('E33').Value = '  +0.68%  '
$ws.Range('E34').Value = '  +0.58%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.57'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.24%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.402.68'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.678'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.27%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('E40').Value = '  +5.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '83.30'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.24%  '
$ws.Range('E42').Value = '  +1.73%  '
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.74'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('E45').Value = '  +2.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0514'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.03'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.984.06'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.21'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.75%  '
$ws.Range('E50').Value = '  +0.62%  '
